$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.423.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.371.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "110.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "314.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.66%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.616"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0922"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.03%  "

$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.983"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.04%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.745.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.48%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.380.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.460.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("E20").Value = "  +1.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.53%  "

$ws.Range("E26").Value = "  -0.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0964"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.47%  "

$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "170.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.70%  "

$ws.Range("E34").Value = "  +4.94%  "

$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.94%  "

$ws.Range("E39").Value = "  +7.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0356"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.235"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.78%  "

$ws.Range("E46").Value = "  -0.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "81.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.647.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.90%  "

